$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text format to preserve exact string formatting ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.370.75"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.390.33"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.98"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.23"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.592"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "679.08"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.65"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.453.57"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.120"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.393.64"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.74"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.31"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.905"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.16"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.25"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.73"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.99"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.75"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "557.46"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.61"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.59"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.674.68"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.78"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0725"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.27"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.70"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0423"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.41"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.63"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"

# --- Coin name / Link / Volume columns ---
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("E3").Value = "  +1.59%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("E9").Value = "  +7.76%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("E12").Value = "  +3.90%  "
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  +1.13%  "
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("E32").Value = "  -2.00%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("E33").Value = "  +8.99%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("E38").Value = "  +1.33%  "
$ws.Range("E39").Value = "  +3.87%  "
$ws.Range("E40").Value = "  +7.55%  "
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  +6.19%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("E51").Value = "  +3.31%  "
